# Helper: assign a text value to a cell without Excel auto-coercing
# numeric-looking strings (e.g. "211.53") into a Number cell. We briefly
# force a Text number-format, set the value, then restore the cell's
# style so the cell's style index ends up unchanged from before the edit.
function Set-Text($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
Set-Text $ws.Range("D2") "89.178.27"
Set-Text $ws.Range("E2") "  -2.63%  "

# Row 3 - Ethereum
Set-Text $ws.Range("D3") "3.129.26"
Set-Text $ws.Range("E3") "  -5.76%  "

# Row 4 - TetherUSD
Set-Text $ws.Range("E4") "  +0.10%  "

# Row 5 - Solana
Set-Text $ws.Range("D5") "211.53"
Set-Text $ws.Range("E5") "  -2.41%  "

# Row 6 - BNB
Set-Text $ws.Range("D6") "624.93"
Set-Text $ws.Range("E6") "  -1.52%  "

# Row 7 - Dogecoin
Set-Text $ws.Range("D7") "0.392"
Set-Text $ws.Range("E7") "  -8.34%  "

# Row 8 - XRP
Set-Text $ws.Range("D8") "0.714"
Set-Text $ws.Range("E8") "  -0.21%  "

# Row 9 - USDC
Set-Text $ws.Range("E9") "  +0.14%  "

# Row 10 - LidoStakedEther
Set-Text $ws.Range("D10") "3.125.87"
Set-Text $ws.Range("E10") "  -5.80%  "

# Row 11 - Cardano
Set-Text $ws.Range("E11") "  -8.79%  "

# Row 12 - TRON
Set-Text $ws.Range("E12") "  -0.21%  "

# Row 13 - ShibaInu
Set-Text $ws.Range("D13") "0.0000249"
Set-Text $ws.Range("E13") "  -9.12%  "

# Row 14 - WrappedBTC
Set-Text $ws.Range("D14") "89.131.50"
Set-Text $ws.Range("E14") "  -2.04%  "

# Row 15 - swap: was WrappedliquidstakedEther2.0, becomes Toncoin
Set-Text $ws.Range("B15") "Toncoin"
Set-Text $ws.Range("C15") "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-Text $ws.Range("D15") "5.24"
Set-Text $ws.Range("E15") "  -4.15%  "

# Row 16 - swap: was Toncoin, becomes WrappedliquidstakedEther2.0
Set-Text $ws.Range("B16") "WrappedliquidstakedEther2.0"
Set-Text $ws.Range("C16") "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-Text $ws.Range("D16") "3.709.83"
Set-Text $ws.Range("E16") "  -4.91%  "

# Row 17 - Avalanche
Set-Text $ws.Range("D17") "31.77"
Set-Text $ws.Range("E17") "  -8.50%  "

# Row 18 - WrappedEther
Set-Text $ws.Range("D18") "3.154.11"
Set-Text $ws.Range("E18") "  -3.86%  "

# Row 19 - SuiNetwork
Set-Text $ws.Range("D19") "3.30"
Set-Text $ws.Range("E19") "  +0.75%  "

# Row 20 - PEPE
Set-Text $ws.Range("E20") "  +14.44%  "

# Row 21 - Chainlink
Set-Text $ws.Range("D21") "13.10"
Set-Text $ws.Range("E21") "  -8.59%  "

# Row 22 - BitcoinCash
Set-Text $ws.Range("D22") "422.32"
Set-Text $ws.Range("E22") "  -3.49%  "

# Row 23 - Uniswap
Set-Text $ws.Range("D23") "8.30"
Set-Text $ws.Range("E23") "  -8.41%  "

# Row 24 - Polkadot
Set-Text $ws.Range("D24") "4.86"
Set-Text $ws.Range("E24") "  -9.88%  "

# Row 25 - NEARProtocol
Set-Text $ws.Range("D25") "5.16"
Set-Text $ws.Range("E25") "  -5.59%  "

# Row 26 - Aptos
Set-Text $ws.Range("D26") "11.38"
Set-Text $ws.Range("E26") "  -6.96%  "

# Row 27 - Litecoin
Set-Text $ws.Range("D27") "79.03"
Set-Text $ws.Range("E27") "  +3.07%  "

# Row 28 - WrappedeETH
Set-Text $ws.Range("D28") "3.351.81"
Set-Text $ws.Range("E28") "  -2.88%  "

# Row 29 - Dai
Set-Text $ws.Range("D29") "0.999"
Set-Text $ws.Range("E29") "  +0.04%  "

# Row 30 - Binance-PegBSC-USD
Set-Text $ws.Range("E30") "  +0.53%  "

# Row 31 - Cronos
Set-Text $ws.Range("D31") "0.155"
Set-Text $ws.Range("E31") "  -12.06%  "

# Row 32 - dogwifhat
Set-Text $ws.Range("D32") "3.89"
Set-Text $ws.Range("E32") "  +7.70%  "

# Row 33 - InternetComputer(DFINITY)
Set-Text $ws.Range("D33") "8.16"
Set-Text $ws.Range("E33") "  -6.55%  "

# Row 34 - Bittensor
Set-Text $ws.Range("D34") "505.34"
Set-Text $ws.Range("E34") "  -11.25%  "

# Row 35 - swap: was PancakeSwap, becomes RenderToken
Set-Text $ws.Range("B35") "RenderToken"
Set-Text $ws.Range("C35") "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-Text $ws.Range("D35") "6.74"
Set-Text $ws.Range("E35") "  -8.80%  "

# Row 36 - swap: was RenderToken, becomes PancakeSwap
Set-Text $ws.Range("B36") "PancakeSwap"
Set-Text $ws.Range("C36") "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-Text $ws.Range("D36") "1.85"
Set-Text $ws.Range("E36") "  -4.35%  "

# Row 37 - Fetch.AI
Set-Text $ws.Range("D37") "1.26"
Set-Text $ws.Range("E37") "  -8.95%  "

# Row 38 - swap: was WhiteBITCoin, becomes EthereumClassic
Set-Text $ws.Range("B38") "EthereumClassic"
Set-Text $ws.Range("C38") "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-Text $ws.Range("D38") "21.85"
Set-Text $ws.Range("E38") "  -4.72%  "

# Row 39 - swap: was EthereumClassic, becomes WhiteBITCoin
Set-Text $ws.Range("B39") "WhiteBITCoin"
Set-Text $ws.Range("C39") "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-Text $ws.Range("D39") "22.24"
Set-Text $ws.Range("E39") "  -0.71%  "

# Row 40 - FirstDigitalUSD
Set-Text $ws.Range("E40") "  +0.06%  "

# Row 41 - Kaspa
Set-Text $ws.Range("E41") "  -7.73%  "

# Row 42 - USDe
Set-Text $ws.Range("E42") "  -0.02%  "

# Row 43 - Stacks
Set-Text $ws.Range("E43") "  -7.49%  "

# Row 44 - PolygonEcosystemToken
Set-Text $ws.Range("D44") "0.362"
Set-Text $ws.Range("E44") "  -8.50%  "

# Row 45 - Monero
Set-Text $ws.Range("D45") "146.74"
Set-Text $ws.Range("E45") "  -1.51%  "

# Row 46 - OKB
Set-Text $ws.Range("D46") "43.54"
Set-Text $ws.Range("E46") "  -2.03%  "

# Row 47 - Aave
Set-Text $ws.Range("D47") "166.27"
Set-Text $ws.Range("E47") "  -11.17%  "

# Row 48 - Stellar
Set-Text $ws.Range("E48") "  -5.05%  "

# Row 49 - Mantle
Set-Text $ws.Range("D49") "0.719"
Set-Text $ws.Range("E49") "  -1.63%  "

# Row 50 - InjectiveProtocol
Set-Text $ws.Range("D50") "24.19"
Set-Text $ws.Range("E50") "  -4.49%  "

# Row 51 - ImmutableX
Set-Text $ws.Range("E51") "  -9.10%  "
